$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.955.82'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = '2.053.29'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.62'
$ws.Range('E7').Value = '  +6.91%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.386'
$ws.Range('E9').Value = '  +2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0813'
$ws.Range('E10').Value = '  +4.06%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = '2.356.39'
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.62'
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('E14').Value = '  +3.70%  '
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  +1.54%  '
$ws.Range('D17').Value = '2.043.72'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').Value = '37.875.01'
$ws.Range('E18').Value = '  +2.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.36'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').Value = '0.0₃0837'
$ws.Range('E21').Value = '  +2.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.48'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.42'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('E25').Value = '  +4.20%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.33'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.28'
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('E28').Value = '  +4.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.03'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.32'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.120'
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('E33').Value = '  +2.95%  '
$ws.Range('E34').Value = '  +10.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0609'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.07'
$ws.Range('E37').Value = '  +13.14%  '
$ws.Range('E38').Value = '  +6.39%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '1.501.93'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0218'
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.10'
$ws.Range('E42').Value = '  +1.99%  '
$ws.Range('E43').Value = '  +3.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.56'
$ws.Range('E44').Value = '  +2.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0918'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.13'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.12'
$ws.Range('E47').Value = '  +15.51%  '
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.09'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('D51').Value = '2.244.03'
$ws.Range('E51').Value = '  +1.54%  '
